# Rename the sole worksheet from "1" to "Bookings", and move the active
# selection from G13 to C20 (matching the saved cursor position recorded
# in the sheet's <selection> element).

$wb = $excel.ActiveWorkbook

# --- Rename sheet "1" -> "Bookings" ------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Bookings"

# --- Update the saved selection/active cell -----------------------------
$ws.Range("C20").Select()

Write-Output ("Sheet1 name: " + $ws.Name)
